# PROD-10305: add/delete row features
# Update the "customers" fixture used by the Cypress tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customers")

# Howard's birthday (E5) is corrected from 21/05/2002 to 12/05/1987.
$ws.Range("E5").Value = "12/05/1987"

# Billy's canDrinkAlcohol flag (C6) becomes a live TRUE() formula instead of
# a hard-coded boolean literal.
$ws.Range("C6").Formula = "=TRUE()"

# Reflect the new selection left behind by the edit (was E7, now E5).
[void]$ws.Range("E5").Select()
